$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2 (A2:AJ2) down to row 3 (A3:AJ3), preserving formatting/values
$ws.Range("A2:AJ2").Copy($ws.Range("A3"))

# Update the specific cells that differ from row 2
$ws.Range("B3").Value = "SYMBOL_2017"
$ws.Range("AC3").Value = "BI001"
$ws.Range("AD3").Value = "PD001"
$ws.Range("AE3").Value = "UM001"
$ws.Range("AF3").Value = "MP001"

# Update the selected cell as in the diff
$ws.Range("F11").Select()
